# Clean up header labels (remove trailing spaces, fix typo "Proovedor",
# collapse "Precio Unitario"/"Precio Total" into single words) and remove
# the trailing space from the "MegaTools" provider name, as part of
# standardizing the report column names for downstream validations.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row (row 1)
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Proveedor"
$ws.Range("C1").Value = "Producto"
$ws.Range("D1").Value = "Cantidad"
$ws.Range("E1").Value = "PrecioUnitario"
$ws.Range("F1").Value = "PrecioTotal"

# Provider name cleanup (remove trailing space on "MegaTools ")
$ws.Range("B2").Value = "MegaTools"
$ws.Range("B3").Value = "MegaTools"
$ws.Range("B4").Value = "MegaTools"
$ws.Range("B5").Value = "MegaTools"

# Update the active selection to match the edited workbook
$ws.Range("E1").Select()
